$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New flight departure rows for "Sunday, Jan 15" added at the bottom of the table
$newRows = @(
    @{ A = 21; B = "Sunday, Jan 15"; C = "10:15 AM"; D = "FR2469"; E = "London";   F = "(STN)"; G = "Ryanair "; H = "B38M"; I = "(EI-HEW)"; J = "10:19 AM"; L = "0 hours, 4 minutes" },
    @{ A = 22; B = "Sunday, Jan 15"; C = "10:30 AM"; D = "FR1979"; E = "Dublin";   F = "(DUB)"; G = "Ryanair "; H = "B38M"; I = "(EI-HMV)"; J = "11:00 AM"; L = "0 hours, 30 minutes" },
    @{ A = 23; B = "Sunday, Jan 15"; C = "10:55 AM"; D = "FR9982"; E = "Alicante"; F = "(ALC)"; G = "Ryanair "; H = "B738"; I = "(EI-DYZ)"; J = "11:06 AM"; L = "0 hours, 11 minutes" }
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    # Columns K and M stay empty in the source data, but the cells still
    # exist in the sheet (same as rows above them) - touch them without
    # altering any formatting so they materialize as blank cells.
    $ws.Cells.Item($r, 11).Borders.LineStyle = -4142
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Borders.LineStyle = -4142
}
